# Apply the diff: add X6/Y6 values to the last existing row, then append a
# new row 7 of scan data (automatic repeater scanner results).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing X6/Y6 cells on the existing last row ---
$ws.Range("X6").Value = -0.45999100000000226
$ws.Range("Y6").Value = "Down"

# --- Append a new row (row 7) with a full scan record ---

# Column A uses the same date/time number format as the rows above it, so
# copy that formatting over before writing the serial date value.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 42648.888831018521

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "Neutral"
$ws.Range("D7").Value = 24
$ws.Range("E7").Value = 28301
$ws.Range("F7").Value = 3299
$ws.Range("G7").Value = 58
$ws.Range("H7").Value = 35
$ws.Range("I7").Value = 82
$ws.Range("J7").Value = 17
$ws.Range("K7").Value = 50150
$ws.Range("L7").Value = 391
$ws.Range("M7").Value = 236
$ws.Range("N7").Value = 121
$ws.Range("O7").Value = 26
$ws.Range("P7").Value = "Bag"
$ws.Range("Q7").Value = 58.594837935340642
$ws.Range("R7").Value = 0

# Columns S and T use the existing percentage number format.
$ws.Range("S6").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = 0.1167

$ws.Range("T6").Copy()
$ws.Range("T7").PasteSpecial(-4122)
$ws.Range("T7").Value = 0.0079000000000000008

$ws.Range("U7").Value = 5.99
$ws.Range("V7").Value = "N/A"
$ws.Range("W7").Value = 0

$excel.CutCopyMode = 0
